$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 24.79234167118218
$ws.Cells.Item(2, 3).Value = 12.19868981780221
$ws.Cells.Item(2, 4).Value = 4.116472192994418
$ws.Cells.Item(2, 5).Value = 9.434954983491286
$ws.Cells.Item(2, 6).Value = 51.05152519456227
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 35.71718101469777
$ws.Cells.Item(2, 10).Value = 9.319292756765989
$ws.Cells.Item(2, 12).Value = 12.18213002633299
$ws.Cells.Item(2, 13).Value = 20.71384453939639
$ws.Cells.Item(2, 14).Value = 22.07692054025606
$ws.Cells.Item(3, 2).Value = 24.52690093257683
$ws.Cells.Item(3, 3).Value = 11.84615480778835
$ws.Cells.Item(3, 4).Value = 4.0880550109427
$ws.Cells.Item(3, 5).Value = 9.422241046720384
$ws.Cells.Item(3, 6).Value = 51.03973662063841
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 35.77477713447134
$ws.Cells.Item(3, 10).Value = 9.325996496786265
$ws.Cells.Item(3, 12).Value = 12.19652173816144
$ws.Cells.Item(3, 13).Value = 20.67533557885137
$ws.Cells.Item(3, 14).Value = 22.14782543832302
$ws.Cells.Item(4, 2).Value = 24.36889176902148
$ws.Cells.Item(4, 3).Value = 11.62822034687833
$ws.Cells.Item(4, 4).Value = 4.070077268847864
$ws.Cells.Item(4, 5).Value = 9.414262293700286
$ws.Cells.Item(4, 6).Value = 51.04436539992822
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 35.81731954102652
$ws.Cells.Item(4, 10).Value = 9.330338538789535
$ws.Cells.Item(4, 12).Value = 12.20699692142302
$ws.Cells.Item(4, 13).Value = 20.65563771639017
$ws.Cells.Item(4, 14).Value = 22.19336377595828
$ws.Cells.Item(5, 2).Value = 24.30582228839021
$ws.Cells.Item(5, 3).Value = 11.53919234269344
$ws.Cells.Item(5, 4).Value = 4.062618182155227
$ws.Cells.Item(5, 5).Value = 9.410967310821276
$ws.Cells.Item(5, 6).Value = 51.04923214866243
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 35.83645571585497
$ws.Cells.Item(5, 10).Value = 9.332164926048851
$ws.Cells.Item(5, 12).Value = 12.21167793815646
$ws.Cells.Item(5, 13).Value = 20.64860851326571
$ws.Cells.Item(5, 14).Value = 22.21242588281552
$ws.Cells.Item(6, 2).Value = 24.29543140194033
$ws.Cells.Item(6, 3).Value = 11.5244007798242
$ws.Cells.Item(6, 4).Value = 4.061371576722818
$ws.Cells.Item(6, 5).Value = 9.410417542507577
$ws.Cells.Item(6, 6).Value = 51.05022014332408
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 35.83974180661274
$ws.Cells.Item(6, 10).Value = 9.33247164183706
$ws.Cells.Item(6, 12).Value = 12.21248012534245
$ws.Cells.Item(6, 13).Value = 20.64750172845293
$ws.Cells.Item(6, 14).Value = 22.21562166299561
$ws.Cells.Item(7, 2).Value = 24.36803575459877
$ws.Cells.Item(7, 3).Value = 11.62702035814402
$ws.Cells.Item(7, 4).Value = 4.06997721135326
$ws.Cells.Item(7, 5).Value = 9.414218033146744
$ws.Cells.Item(7, 6).Value = 51.04441897298004
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 35.81757033817435
$ws.Cells.Item(7, 10).Value = 9.330362939204674
$ws.Cells.Item(7, 12).Value = 12.20705838160088
$ws.Cells.Item(7, 13).Value = 20.65553887141188
$ws.Cells.Item(7, 14).Value = 22.19361880830841
$ws.Cells.Item(8, 2).Value = 24.69982855719711
$ws.Cells.Item(8, 3).Value = 12.07753171101105
$ws.Cells.Item(8, 4).Value = 4.106783915897962
$ws.Cells.Item(8, 5).Value = 9.430606842699856
$ws.Cells.Item(8, 6).Value = 51.04499585664034
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 35.73554730859874
$ws.Cells.Item(8, 10).Value = 9.321557413917327
$ws.Cells.Item(8, 12).Value = 12.18675235847364
$ws.Cells.Item(8, 13).Value = 20.69975094396984
$ws.Cells.Item(8, 14).Value = 22.10095374863339
$ws.Cells.Item(9, 2).Value = 25.38672562499629
$ws.Cells.Item(9, 3).Value = 12.94295895982826
$ws.Cells.Item(9, 4).Value = 4.174756801362238
$ws.Cells.Item(9, 5).Value = 9.461389392393007
$ws.Cells.Item(9, 6).Value = 51.14037594266006
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 35.63187536417949
$ws.Cells.Item(9, 10).Value = 9.306075025469001
$ws.Cells.Item(9, 12).Value = 12.15992201977785
$ws.Cells.Item(9, 13).Value = 20.81749914017265
$ws.Cells.Item(9, 14).Value = 21.93506412932303
$ws.Cells.Item(10, 2).Value = 25.90907205686634
$ws.Cells.Item(10, 3).Value = 13.56001811907462
$ws.Cells.Item(10, 4).Value = 4.222131006246007
$ws.Cells.Item(10, 5).Value = 9.483196077570259
$ws.Cells.Item(10, 6).Value = 51.267928029023
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 35.59084323718587
$ws.Cells.Item(10, 10).Value = 9.29577815138194
$ws.Cells.Item(10, 12).Value = 12.14811119368459
$ws.Cells.Item(10, 13).Value = 20.92253340934198
$ws.Cells.Item(10, 14).Value = 21.82274419021705
$ws.Cells.Item(11, 2).Value = 26.1495800097273
$ws.Cells.Item(11, 3).Value = 13.83521420947759
$ws.Cells.Item(11, 4).Value = 4.243125432202477
$ws.Cells.Item(11, 5).Value = 9.492944529659797
$ws.Cells.Item(11, 6).Value = 51.33840090837135
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 35.57985725244958
$ws.Cells.Item(11, 10).Value = 9.291325737853304
$ws.Cells.Item(11, 12).Value = 12.14444908991992
$ws.Cells.Item(11, 13).Value = 20.97423735152845
$ws.Cells.Item(11, 14).Value = 21.77370371715278
$ws.Cells.Item(12, 2).Value = 26.24098889567962
$ws.Cells.Item(12, 3).Value = 13.93852551327057
$ws.Cells.Item(12, 4).Value = 4.250995370109025
$ws.Cells.Item(12, 5).Value = 9.496611614127584
$ws.Cells.Item(12, 6).Value = 51.3668709037321
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 35.57680480177344
$ws.Cells.Item(12, 10).Value = 9.289672876170894
$ws.Cells.Item(12, 12).Value = 12.14330781508883
$ws.Cells.Item(12, 13).Value = 20.99437066254046
$ws.Cells.Item(12, 14).Value = 21.75542742607267
$ws.Cells.Item(13, 2).Value = 26.22128880977843
$ws.Cells.Item(13, 3).Value = 13.91631718488012
$ws.Cells.Item(13, 4).Value = 4.249304007951938
$ws.Cells.Item(13, 5).Value = 9.495822925471591
$ws.Cells.Item(13, 6).Value = 51.36066018860544
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 35.57741288795896
$ws.Cells.Item(13, 10).Value = 9.290027376535583
$ws.Cells.Item(13, 12).Value = 12.14354269989303
$ws.Cells.Item(13, 13).Value = 20.99001011472367
$ws.Cells.Item(13, 14).Value = 21.75935048269486
$ws.Cells.Item(14, 2).Value = 26.15709405171163
$ws.Cells.Item(14, 3).Value = 13.84373229952807
$ws.Cells.Item(14, 4).Value = 4.243774505067347
$ws.Cells.Item(14, 5).Value = 9.493246711512256
$ws.Cells.Item(14, 6).Value = 51.34070745429958
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 35.57958390798598
$ws.Cells.Item(14, 10).Value = 9.291189091966999
$ws.Cells.Item(14, 12).Value = 12.1443502798408
$ws.Cells.Item(14, 13).Value = 20.97588268501406
$ws.Cells.Item(14, 14).Value = 21.77219422458529
$ws.Cells.Item(15, 2).Value = 26.1178139530698
$ws.Cells.Item(15, 3).Value = 13.79915180377466
$ws.Cells.Item(15, 4).Value = 4.24037707535502
$ws.Cells.Item(15, 5).Value = 9.491665529449582
$ws.Cells.Item(15, 6).Value = 51.32871785420948
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 35.58105806741819
$ws.Cells.Item(15, 10).Value = 9.291904991719473
$ws.Cells.Item(15, 12).Value = 12.14487689890062
$ws.Cells.Item(15, 13).Value = 20.96730108383412
$ws.Cells.Item(15, 14).Value = 21.78009967629988
$ws.Cells.Item(16, 2).Value = 25.89340544380607
$ws.Cells.Item(16, 3).Value = 13.54191316125435
$ws.Cells.Item(16, 4).Value = 4.220747759447923
$ws.Cells.Item(16, 5).Value = 9.482555560874919
$ws.Cells.Item(16, 6).Value = 51.26357243607067
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 35.59171603227371
$ws.Cells.Item(16, 10).Value = 9.296073776727823
$ws.Cells.Item(16, 12).Value = 12.14838490210634
$ws.Cells.Item(16, 13).Value = 20.91923254204862
$ws.Cells.Item(16, 14).Value = 21.82599034189374
$ws.Cells.Item(17, 2).Value = 25.75641992205315
$ws.Cells.Item(17, 3).Value = 13.38261528453143
$ws.Cells.Item(17, 4).Value = 4.208563178644017
$ws.Cells.Item(17, 5).Value = 9.476923118327123
$ws.Cells.Item(17, 6).Value = 51.22679213703682
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 35.6002238669371
$ws.Cells.Item(17, 10).Value = 9.298690427746708
$ws.Cells.Item(17, 12).Value = 12.15097475550454
$ws.Cells.Item(17, 13).Value = 20.89074200213276
$ws.Cells.Item(17, 14).Value = 21.85466814171421
$ws.Cells.Item(18, 2).Value = 25.67790637175003
$ws.Cells.Item(18, 3).Value = 13.29047876752263
$ws.Cells.Item(18, 4).Value = 4.201502576353875
$ws.Cells.Item(18, 5).Value = 9.473667329097433
$ws.Cells.Item(18, 6).Value = 51.20680934567975
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 35.60584003177416
$ws.Cells.Item(18, 10).Value = 9.300217272223628
$ws.Cells.Item(18, 12).Value = 12.15262542457683
$ws.Cells.Item(18, 13).Value = 20.87472492445568
$ws.Cells.Item(18, 14).Value = 21.87135629426986
$ws.Cells.Item(19, 2).Value = 25.65137309080994
$ws.Cells.Item(19, 3).Value = 13.2591983324632
$ws.Cells.Item(19, 4).Value = 4.199102991963056
$ws.Cells.Item(19, 5).Value = 9.472562193019124
$ws.Cells.Item(19, 6).Value = 51.200245024591
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 35.60786557935918
$ws.Cells.Item(19, 10).Value = 9.300737986903169
$ws.Cells.Item(19, 12).Value = 12.153211988474
$ws.Cells.Item(19, 13).Value = 20.86936564623058
$ws.Cells.Item(19, 14).Value = 21.87703987497222
$ws.Cells.Item(20, 2).Value = 25.77097420122354
$ws.Cells.Item(20, 3).Value = 13.39962667813191
$ws.Cells.Item(20, 4).Value = 4.209865667849551
$ws.Cells.Item(20, 5).Value = 9.477524372024613
$ws.Cells.Item(20, 6).Value = 51.23058618038407
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 35.59924337275117
$ws.Cells.Item(20, 10).Value = 9.298409623965377
$ws.Cells.Item(20, 12).Value = 12.15068239626282
$ws.Cells.Item(20, 13).Value = 20.89373665312884
$ws.Cells.Item(20, 14).Value = 21.85159532913422
$ws.Cells.Item(21, 2).Value = 26.17594119980112
$ws.Cells.Item(21, 3).Value = 13.86507745832292
$ws.Cells.Item(21, 4).Value = 4.245400832118766
$ws.Cells.Item(21, 5).Value = 9.494004069615372
$ws.Cells.Item(21, 6).Value = 51.34651971668674
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 35.57891613968624
$ws.Cells.Item(21, 10).Value = 9.290846968894511
$ws.Cells.Item(21, 12).Value = 12.1441064160562
$ws.Cells.Item(21, 13).Value = 20.98001729739505
$ws.Cells.Item(21, 14).Value = 21.76841372889868
$ws.Cells.Item(22, 2).Value = 26.44251841377987
$ws.Cells.Item(22, 3).Value = 14.16399078580409
$ws.Cells.Item(22, 4).Value = 4.268157588183342
$ws.Cells.Item(22, 5).Value = 9.504632135079873
$ws.Cells.Item(22, 6).Value = 51.43267980058817
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 35.57208886428297
$ws.Cells.Item(22, 10).Value = 9.286097614192535
$ws.Cells.Item(22, 12).Value = 12.14123922095858
$ws.Cells.Item(22, 13).Value = 21.03963156551475
$ws.Cells.Item(22, 14).Value = 21.71576446347416
$ws.Cells.Item(23, 2).Value = 26.30009342809098
$ws.Cells.Item(23, 3).Value = 14.00497167792417
$ws.Cells.Item(23, 4).Value = 4.256054721202577
$ws.Cells.Item(23, 5).Value = 9.498972690635069
$ws.Cells.Item(23, 6).Value = 51.38574648995844
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 35.57514084371186
$ws.Cells.Item(23, 10).Value = 9.28861479700336
$ws.Cells.Item(23, 12).Value = 12.14263878373348
$ws.Cells.Item(23, 13).Value = 21.00752272035562
$ws.Cells.Item(23, 14).Value = 21.74370784655208
$ws.Cells.Item(24, 2).Value = 25.76439345387807
$ws.Cells.Item(24, 3).Value = 13.39193754458104
$ws.Cells.Item(24, 4).Value = 4.20927698543958
$ws.Cells.Item(24, 5).Value = 9.477252600152974
$ws.Cells.Item(24, 6).Value = 51.22886727116575
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 35.5996843963878
$ws.Cells.Item(24, 10).Value = 9.298536505269141
$ws.Cells.Item(24, 12).Value = 12.1508140680765
$ws.Cells.Item(24, 13).Value = 20.89238164145952
$ws.Cells.Item(24, 14).Value = 21.85298392178031
$ws.Cells.Item(25, 2).Value = 25.19747164687396
$ws.Cells.Item(25, 3).Value = 12.71160967893754
$ws.Cells.Item(25, 4).Value = 4.156818706975544
$ws.Cells.Item(25, 5).Value = 9.453205252937503
$ws.Cells.Item(25, 6).Value = 51.10447357438157
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 35.65376970399458
$ws.Cells.Item(25, 10).Value = 9.310073369184099
$ws.Cells.Item(25, 12).Value = 12.165791269619
$ws.Cells.Item(25, 13).Value = 20.78235978823934
$ws.Cells.Item(25, 14).Value = 21.97825661289538

Write-Output "Updated 264 cells (rows 2-25, cols B-N excl. G/K)"
